$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: "Trening" header + "Gra" for every data row ---
$ws.Range("F1").Value = "Trening"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats - match header style (bold, border, centered)
$excel.CutCopyMode = $false

$ws.Range("F2:F7").Value = "Gra"

# --- Column A: convert text timestamps to real Excel date-time serials ---
# Apply the custom format once to a single cell first (registers numFmtId 164
# with the lowercase code), then re-apply the final uppercase code to that
# same cell (registers numFmtId 165 and switches the cell's style in place),
# then extend that style to the rest of the column.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A3:A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A2").Value = 45687.521828125
$ws.Range("A3").Value = 45687.5296568287
$ws.Range("A4").Value = 45687.52990682871
$ws.Range("A5").Value = 45687.51643576389
$ws.Range("A6").Value = 45687.52965567129
$ws.Range("A7").Value = 45687.5299056713
